$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.0089430809021
$ws.Range("B1").Value = 4.266225814819336
$ws.Range("C1").Value = 8.006000518798828
$ws.Range("D1").Value = 8.246187210083008
$ws.Range("E1").Value = 5.573213577270508
